$wb = $excel.ActiveWorkbook

# ---- ALC (sheet index 1) ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 1912.1538
$ws.Cells.Item(19, 9).Value = 3644.8333
$ws.Cells.Item(19, 10).Value = 427
$ws.Cells.Item(19, 11).Value = 3644.8333
$ws.Cells.Item(19, 12).Value = 427
$ws.Cells.Item(19, 13).Value = -3469.8333
$ws.Cells.Item(19, 14).Value = -777

$ws.Cells.Item(62, 8).Value = 6891.5264
$ws.Cells.Item(62, 9).Value = 6696
$ws.Cells.Item(62, 11).Value = 6696
$ws.Cells.Item(62, 13).Value = -6072

$ws.Cells.Item(65, 8).Value = 6891.5264
$ws.Cells.Item(65, 9).Value = 6696
$ws.Cells.Item(65, 11).Value = 33480
$ws.Cells.Item(65, 13).Value = -30360

$ws.Cells.Item(100, 8).Value = 29257.406
$ws.Cells.Item(100, 10).Value = 3665.2727
$ws.Cells.Item(100, 12).Value = 3665.2727
$ws.Cells.Item(100, 14).Value = -4747.2727

$ws.Cells.Item(137, 8).Value = 9796.796
$ws.Cells.Item(137, 9).Value = 5038.0454
$ws.Cells.Item(137, 11).Value = 15114.1362
$ws.Cells.Item(137, 13).Value = -12564.1362

$ws.Cells.Item(141, 8).Value = 1521.1786
$ws.Cells.Item(141, 9).Value = 1538.4445
$ws.Cells.Item(141, 11).Value = 4615.333500000001
$ws.Cells.Item(141, 13).Value = 564.6664999999994

# ---- ARM (sheet index 2) ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 9149.058999999999
$ws.Cells.Item(61, 9).Value = 6563.76
$ws.Cells.Item(61, 10).Value = 16330.444
$ws.Cells.Item(61, 11).Value = 6563.76
$ws.Cells.Item(61, 12).Value = 16330.444
$ws.Cells.Item(61, 13).Value = -6351.76
$ws.Cells.Item(61, 14).Value = -16754.444

$ws.Cells.Item(74, 8).Value = 4490.6523
$ws.Cells.Item(74, 9).Value = 3924.4666
$ws.Cells.Item(74, 10).Value = 5552.25
$ws.Cells.Item(74, 11).Value = 3924.4666
$ws.Cells.Item(74, 12).Value = 5552.25
$ws.Cells.Item(74, 13).Value = -3050.4666
$ws.Cells.Item(74, 14).Value = -7300.25

$ws.Cells.Item(77, 8).Value = 4490.6523
$ws.Cells.Item(77, 9).Value = 3924.4666
$ws.Cells.Item(77, 10).Value = 5552.25
$ws.Cells.Item(77, 11).Value = 19622.333
$ws.Cells.Item(77, 12).Value = 27761.25
$ws.Cells.Item(77, 13).Value = -15254.333
$ws.Cells.Item(77, 14).Value = -36497.25

$ws.Cells.Item(110, 8).Value = 374.2
$ws.Cells.Item(110, 9).Value = 318
$ws.Cells.Item(110, 11).Value = 318
$ws.Cells.Item(110, 13).Value = 1727

$ws.Cells.Item(132, 8).Value = 3211.7837
$ws.Cells.Item(132, 9).Value = 2935.875
$ws.Cells.Item(132, 11).Value = 8807.625
$ws.Cells.Item(132, 13).Value = -6277.625

$ws.Cells.Item(136, 8).Value = 9149.058999999999
$ws.Cells.Item(136, 9).Value = 6563.76
$ws.Cells.Item(136, 10).Value = 16330.444
$ws.Cells.Item(136, 11).Value = 19691.28
$ws.Cells.Item(136, 12).Value = 48991.33199999999
$ws.Cells.Item(136, 13).Value = -17141.28
$ws.Cells.Item(136, 14).Value = -54091.33199999999

# ---- BSM (sheet index 3) ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 15280.038
$ws.Cells.Item(20, 9).Value = 22271.75
$ws.Cells.Item(20, 10).Value = 4093.3
$ws.Cells.Item(20, 11).Value = 22271.75
$ws.Cells.Item(20, 12).Value = 4093.3
$ws.Cells.Item(20, 13).Value = -22024.75
$ws.Cells.Item(20, 14).Value = -4587.3

$ws.Cells.Item(94, 8).Value = 1084.8975
$ws.Cells.Item(94, 9).Value = 1088.0571
$ws.Cells.Item(94, 11).Value = 1088.0571
$ws.Cells.Item(94, 13).Value = -637.0571

$ws.Cells.Item(134, 8).Value = 5469.449
$ws.Cells.Item(134, 9).Value = 4046.6897
$ws.Cells.Item(134, 11).Value = 12140.0691
$ws.Cells.Item(134, 13).Value = -9605.069100000001

$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 13).ClearContents()

$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

# ---- CRP (sheet index 4) ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2700.3901
$ws.Cells.Item(31, 9).Value = 1916.6666
$ws.Cells.Item(31, 10).Value = 3313.739
$ws.Cells.Item(31, 11).Value = 1916.6666
$ws.Cells.Item(31, 12).Value = 3313.739
$ws.Cells.Item(31, 13).Value = -1621.6666
$ws.Cells.Item(31, 14).Value = -3903.739

$ws.Cells.Item(34, 8).Value = 2700.3901
$ws.Cells.Item(34, 9).Value = 1916.6666
$ws.Cells.Item(34, 10).Value = 3313.739
$ws.Cells.Item(34, 11).Value = 1916.6666
$ws.Cells.Item(34, 12).Value = 3313.739
$ws.Cells.Item(34, 13).Value = -1714.6666
$ws.Cells.Item(34, 14).Value = -3717.739

$ws.Cells.Item(58, 8).Value = 2986.7568
$ws.Cells.Item(58, 9).Value = 1350.0333
$ws.Cells.Item(58, 11).Value = 1350.0333
$ws.Cells.Item(58, 13).Value = -1147.0333

$ws.Cells.Item(132, 8).Value = 22985.443
$ws.Cells.Item(132, 9).Value = 14256.857
$ws.Cells.Item(132, 10).Value = 42280.21
$ws.Cells.Item(132, 11).Value = 42770.571
$ws.Cells.Item(132, 12).Value = 126840.63
$ws.Cells.Item(132, 13).Value = -40240.571
$ws.Cells.Item(132, 14).Value = -131900.63

$ws.Cells.Item(134, 8).Value = 3664
$ws.Cells.Item(134, 9).Value = 2296.0967
$ws.Cells.Item(134, 10).Value = 7519
$ws.Cells.Item(134, 11).Value = 6888.2901
$ws.Cells.Item(134, 12).Value = 22557
$ws.Cells.Item(134, 13).Value = -4353.2901
$ws.Cells.Item(134, 14).Value = -27627

$ws.Cells.Item(136, 8).Value = 2986.7568
$ws.Cells.Item(136, 9).Value = 1350.0333
$ws.Cells.Item(136, 11).Value = 4050.0999
$ws.Cells.Item(136, 13).Value = -1500.0999

# ---- CUL (sheet index 5) ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(107, 8).Value = 584.9666999999999
$ws.Cells.Item(107, 9).Value = 569.1053000000001
$ws.Cells.Item(107, 10).Value = 612.36365
$ws.Cells.Item(107, 11).Value = 1707.3159
$ws.Cells.Item(107, 12).Value = 1837.09095
$ws.Cells.Item(107, 13).Value = 212.6840999999999
$ws.Cells.Item(107, 14).Value = -5677.09095

$ws.Cells.Item(131, 8).Value = 1356.0975
$ws.Cells.Item(131, 10).Value = 1368.5135
$ws.Cells.Item(131, 12).Value = 4105.5405
$ws.Cells.Item(131, 14).Value = -14185.5405

# ---- GSM (sheet index 6) ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(10, 8).Value = 15393.8
$ws.Cells.Item(10, 9).Value = 15663
$ws.Cells.Item(10, 10).Value = 14990
$ws.Cells.Item(10, 11).Value = 15663
$ws.Cells.Item(10, 12).Value = 14990
$ws.Cells.Item(10, 13).Value = -15494
$ws.Cells.Item(10, 14).Value = -15328

$ws.Cells.Item(122, 8).Value = 3416.6667
$ws.Cells.Item(122, 9).Value = 3455.3635
$ws.Cells.Item(122, 11).Value = 10366.0905
$ws.Cells.Item(122, 13).Value = -7916.0905

$ws.Cells.Item(132, 8).Value = 3999.9038
$ws.Cells.Item(132, 9).Value = 1631.5385
$ws.Cells.Item(132, 11).Value = 4894.6155
$ws.Cells.Item(132, 13).Value = -2364.6155

# ---- LTW (sheet index 7) ----
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(14, 8).Value = 14493
$ws.Cells.Item(14, 9).Value = 14493
$ws.Cells.Item(14, 11).Value = 14493
$ws.Cells.Item(14, 13).Value = -14321

$ws.Cells.Item(22, 9).Value = 1569.1818
$ws.Cells.Item(22, 10).Value = 1333.3334
$ws.Cells.Item(22, 11).Value = 1569.1818
$ws.Cells.Item(22, 12).Value = 1333.3334
$ws.Cells.Item(22, 13).Value = -1274.1818
$ws.Cells.Item(22, 14).Value = -1923.3334

$ws.Cells.Item(27, 9).Value = 1569.1818
$ws.Cells.Item(27, 10).Value = 1333.3334
$ws.Cells.Item(27, 11).Value = 1569.1818
$ws.Cells.Item(27, 12).Value = 1333.3334
$ws.Cells.Item(27, 13).Value = -1462.1818
$ws.Cells.Item(27, 14).Value = -1547.3334

$ws.Cells.Item(46, 8).Value = 1361.4
$ws.Cells.Item(46, 9).Value = 1127.8572
$ws.Cells.Item(46, 11).Value = 1127.8572
$ws.Cells.Item(46, 13).Value = -939.8571999999999

$ws.Cells.Item(122, 8).Value = 3309.5557
$ws.Cells.Item(122, 10).Value = 3988
$ws.Cells.Item(122, 12).Value = 11964
$ws.Cells.Item(122, 14).Value = -16864

$ws.Cells.Item(136, 8).Value = 5914.643
$ws.Cells.Item(136, 9).Value = 5345.375
$ws.Cells.Item(136, 11).Value = 16036.125
$ws.Cells.Item(136, 13).Value = -13486.125

# ---- WVR (sheet index 8) ----
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 2392.0527
$ws.Cells.Item(81, 9).Value = 1844.9
$ws.Cells.Item(81, 11).Value = 3689.8
$ws.Cells.Item(81, 13).Value = -2628.8

$ws.Cells.Item(84, 8).Value = 2392.0527
$ws.Cells.Item(84, 9).Value = 1844.9
$ws.Cells.Item(84, 11).Value = 18449
$ws.Cells.Item(84, 13).Value = -13145

$ws.Cells.Item(122, 8).Value = 4838.913
$ws.Cells.Item(122, 9).Value = 2989.2632
$ws.Cells.Item(122, 11).Value = 8967.7896
$ws.Cells.Item(122, 13).Value = -6517.7896

$ws.Cells.Item(132, 8).Value = 12732.478
$ws.Cells.Item(132, 9).Value = 6993.754
$ws.Cells.Item(132, 10).Value = 25697.74
$ws.Cells.Item(132, 11).Value = 20981.262
$ws.Cells.Item(132, 12).Value = 77093.22
$ws.Cells.Item(132, 13).Value = -18451.262
$ws.Cells.Item(132, 14).Value = -82153.22

$ws.Cells.Item(135, 8).Value = 72060
$ws.Cells.Item(135, 10).Value = 72060
$ws.Cells.Item(135, 12).Value = 72060
$ws.Cells.Item(135, 14).Value = -82200

$ws.Cells.Item(136, 8).Value = 1502.6
$ws.Cells.Item(136, 9).Value = 1467.0714
$ws.Cells.Item(136, 11).Value = 4401.2142
$ws.Cells.Item(136, 13).Value = -1851.2142
